$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 640733.9399999999
$ws.Range("J17").Value = 640733.9399999999
$ws.Range("L17").Value = 1922201.82
$ws.Range("N17").Value = -1922537.82

$ws.Range("H40").Value = 2599
$ws.Range("J40").Value = 2599
$ws.Range("L40").Value = 2599
$ws.Range("N40").Value = -2949

$ws.Range("H86").Value = 8699.200000000001
$ws.Range("I86").Value = 7834
$ws.Range("J86").Value = 9997
$ws.Range("K86").Value = 7834
$ws.Range("L86").Value = 9997
$ws.Range("M86").Value = -6711
$ws.Range("N86").Value = -12243

$ws.Range("H89").Value = 8699.200000000001
$ws.Range("I89").Value = 7834
$ws.Range("J89").Value = 9997
$ws.Range("K89").Value = 39170
$ws.Range("L89").Value = 49985
$ws.Range("M89").Value = -33554
$ws.Range("N89").Value = -61217

$ws.Range("H133").Value = 77987.5
$ws.Range("J133").Value = 77987.5
$ws.Range("L133").Value = 77987.5
$ws.Range("N133").Value = -88107.5

$ws.Range("H138").Value = 28160.62
$ws.Range("I138").Value = 69224.92999999999
$ws.Range("J138").Value = 5347.1113
$ws.Range("K138").Value = 207674.79
$ws.Range("L138").Value = 16041.3339
$ws.Range("M138").Value = -202534.79
$ws.Range("N138").Value = -26321.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5158.4653
$ws.Range("I32").Value = 4674.795
$ws.Range("J32").Value = 9874.25
$ws.Range("K32").Value = 4674.795
$ws.Range("L32").Value = 9874.25
$ws.Range("M32").Value = -4387.795
$ws.Range("N32").Value = -10448.25

$ws.Range("H55").Value = 7991.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H122").Value = 1723.5555
$ws.Range("I122").Value = 1503
$ws.Range("K122").Value = 4509
$ws.Range("M122").Value = -2059

$ws.Range("H138").Value = 80195
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1748.921
$ws.Range("I107").Value = 1313.3928
$ws.Range("K107").Value = 1313.3928
$ws.Range("M107").Value = 606.6071999999999

$ws.Range("H132").Value = 100867.8
$ws.Range("J132").Value = 100867.8
$ws.Range("L132").Value = 100867.8
$ws.Range("N132").Value = -110987.8

$ws.Range("H134").Value = 2011.5151
$ws.Range("I134").Value = 2027.5
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 6082.5
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -3547.5
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 119.71429
$ws.Range("I8").Value = 119.71429
$ws.Range("K8").Value = 359.14287
$ws.Range("M8").Value = -220.14287

$ws.Range("H131").Value = 33520.094
$ws.Range("I131").Value = 112010.664
$ws.Range("J131").Value = 2806.3914
$ws.Range("K131").Value = 336031.992
$ws.Range("L131").Value = 8419.174199999999
$ws.Range("M131").Value = -330991.992
$ws.Range("N131").Value = -18499.1742

$ws.Range("H132").Value = 1127.8667
$ws.Range("J132").Value = 1300
$ws.Range("L132").Value = 11700
$ws.Range("N132").Value = -16760

$ws.Range("H139").Value = 3736.75
$ws.Range("I139").Value = 2966.6667
$ws.Range("J139").Value = 4198.8
$ws.Range("K139").Value = 8900.000100000001
$ws.Range("L139").Value = 12596.4
$ws.Range("M139").Value = -3760.000100000001
$ws.Range("N139").Value = -22876.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 39000
$ws.Range("I62").Value = 39000
$ws.Range("K62").Value = 39000
$ws.Range("M62").Value = -38314

$ws.Range("H65").Value = 39000
$ws.Range("I65").Value = 39000
$ws.Range("K65").Value = 117000
$ws.Range("M65").Value = -113568

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H113").Value = 2140.25
$ws.Range("I113").Value = 2018.7693
$ws.Range("K113").Value = 2018.7693
$ws.Range("M113").Value = 151.2307000000001

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 29998.834
$ws.Range("J120").Value = 29998.834
$ws.Range("L120").Value = 29998.834
$ws.Range("N120").Value = -39674.834

$ws.Range("H122").Value = 1417.7222
$ws.Range("I122").Value = 1212.8823
$ws.Range("K122").Value = 3638.6469
$ws.Range("M122").Value = -1188.6469

$ws.Range("H126").Value = 3786.4856
$ws.Range("I126").Value = 3111.6
$ws.Range("K126").Value = 9334.799999999999
$ws.Range("M126").Value = -6864.799999999999

$ws.Range("H141").Value = 78995.336
$ws.Range("J141").Value = 78995.336
$ws.Range("L141").Value = 78995.336
$ws.Range("N141").Value = -89355.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8748.951999999999
$ws.Range("I7").Value = 8985.684999999999
$ws.Range("J7").Value = 6500
$ws.Range("K7").Value = 8985.684999999999
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = -8873.684999999999
$ws.Range("N7").Value = -6724

$ws.Range("H43").Value = 14755.875
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 15578.143
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 15578.143
$ws.Range("M43").Value = -8807
$ws.Range("N43").Value = -15964.143

$ws.Range("H55").Value = 1278
$ws.Range("I55").Value = 900
$ws.Range("K55").Value = 900
$ws.Range("M55").Value = -727

$ws.Range("H62").Value = 39999.668
$ws.Range("J62").Value = 41000
$ws.Range("L62").Value = 41000
$ws.Range("N62").Value = -42248

$ws.Range("H65").Value = 39999.668
$ws.Range("J65").Value = 41000
$ws.Range("L65").Value = 123000
$ws.Range("N65").Value = -129240

$ws.Range("H126").Value = 8748.951999999999
$ws.Range("I126").Value = 8985.684999999999
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 26957.055
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -24487.055
$ws.Range("N126").Value = -24440

$ws.Range("H132").Value = 4995.0586
$ws.Range("I132").Value = 4922.2144
$ws.Range("J132").Value = 5335
$ws.Range("K132").Value = 14766.6432
$ws.Range("L132").Value = 16005
$ws.Range("M132").Value = -12236.6432
$ws.Range("N132").Value = -21065

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4979.6665
$ws.Range("I62").Value = 3487.5
$ws.Range("J62").Value = 5725.75
$ws.Range("K62").Value = 3487.5
$ws.Range("L62").Value = 5725.75
$ws.Range("M62").Value = -2863.5
$ws.Range("N62").Value = -6973.75

$ws.Range("H65").Value = 4979.6665
$ws.Range("I65").Value = 3487.5
$ws.Range("J65").Value = 5725.75
$ws.Range("K65").Value = 17437.5
$ws.Range("L65").Value = 28628.75
$ws.Range("M65").Value = -14317.5
$ws.Range("N65").Value = -34868.75

$ws.Range("H100").Value = 1616.3334
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 1674.6428
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 3349.2856
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -4431.2856

$ws.Range("H122").Value = 3077.8206
$ws.Range("I122").Value = 2833.516
$ws.Range("J122").Value = 4024.5
$ws.Range("K122").Value = 8500.548000000001
$ws.Range("L122").Value = 12073.5
$ws.Range("M122").Value = -6050.548000000001
$ws.Range("N122").Value = -16973.5

$ws.Range("H126").Value = 3486
$ws.Range("I126").Value = 2950.7778
$ws.Range("J126").Value = 5894.5
$ws.Range("K126").Value = 8852.3334
$ws.Range("L126").Value = 17683.5
$ws.Range("M126").Value = -6382.3334
$ws.Range("N126").Value = -22623.5
